$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-18 Thursday" "2025-09-19 Friday"

Replace-Text "722÷7=" "181÷3="
Replace-Text "730÷5=" "716÷5="
Replace-Text "503÷6=" "245÷6="
Replace-Text "655÷3=" "347÷7="
Replace-Text "556÷4=" "670÷2="
Replace-Text "660÷9=" "607÷2="
Replace-Text "770÷8=" "635÷7="
Replace-Text "293÷7=" "547÷7="
Replace-Text "330÷7=" "499÷7="
Replace-Text "142÷2=" "669÷8="
Replace-Text "618÷4=" "820÷8="
Replace-Text "726÷7=" "826÷9="
Replace-Text "749÷6=" "426÷3="
Replace-Text "878÷4=" "940÷3="
Replace-Text "947÷8=" "502÷8="
Replace-Text "201÷4=" "586÷9="
Replace-Text "568÷9=" "565÷3="
Replace-Text "698÷2=" "546÷3="
Replace-Text "996÷8=" "763÷2="
Replace-Text "120÷3=" "314÷2="
Replace-Text "144÷2=" "871÷2="
Replace-Text "430÷6=" "538÷5="
Replace-Text "890÷6=" "733÷4="
Replace-Text "814÷7=" "259÷4="
Replace-Text "123÷6=" "301÷2="
